$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 08.10.2024"

# Row 6
$ws.Range("B6").Value = "09.10."
$ws.Range("C6").Value = "10.10."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 97418195"
$ws.Range("E6").Value = "40,57-"

# Row 7
$ws.Range("B7").Value = "13.10."
$ws.Range("C7").Value = "14.10."
$ws.Range("D7").Value = "KARTENZ./13.10 ALDI SUED RO"
$ws.Range("E7").Value = "90,25-"

# Row 8
$ws.Range("B8").Value = "14.10."
$ws.Range("C8").Value = "15.10."
$ws.Range("D8").Value = "KARTENZ./14.10 LIDL RO"
$ws.Range("E8").Value = "112,58-"

# Row 9
$ws.Range("B9").Value = "16.10."
$ws.Range("C9").Value = "17.10."
$ws.Range("D9").Value = "BURGER KING Siegen"
$ws.Range("E9").Value = "21,45-"

# Row 10
$ws.Range("B10").Value = "17.10."
$ws.Range("C10").Value = "18.10."
$ws.Range("D10").Value = "PAYPAL BABYXU"
$ws.Range("E10").Value = "66,27-"

# Row 11 - previously blank, now a new transaction row
$ws.Range("B11").Value = "18.10."
$ws.Range("C11").Value = "19.10."
$ws.Range("D11").Value = "KARTENZAHLUNG JET TANKSTELLE"
# E11's format needs to change from its old (blank-row) style to the
# transaction-amount style used by E6:E10 - copy format from E10 first.
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "64,88-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 20.10.2024"
$ws.Range("E12").Value = "396,00-"

# Next billing date note
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.10.2024"
